{"js": "// Replace each three-digit-by-one-digit multiplication expression in the\n// document with its new value, per the commit's regenerated numbers.\n// Every \"old\" string below is unique in the document, so searching with\n// matchCase true and replacing the single hit found is unambiguous.\nconst replacements = [\n  [\"195\u00d77=\", \"761\u00d75=\"],\n  [\"552\u00d78=\", \"230\u00d73=\"],\n  [\"709\u00d79=\", \"367\u00d75=\"],\n  [\"248\u00d76=\", \"844\u00d79=\"],\n  [\"234\u00d77=\", \"418\u00d72=\"],\n  [\"305\u00d79=\", \"181\u00d75=\"],\n  [\"260\u00d74=\", \"261\u00d73=\"],\n  [\"409\u00d73=\", \"423\u00d76=\"],\n  [\"314\u00d77=\", \"721\u00d73=\"],\n  [\"149\u00d74=\", \"935\u00d72=\"],\n  [\"243\u00d75=\", \"698\u00d72=\"],\n  [\"386\u00d72=\", \"484\u00d78=\"],\n  [\"182\u00d77=\", \"513\u00d72=\"],\n  [\"976\u00d75=\", \"814\u00d74=\"],\n  [\"778\u00d75=\", \"620\u00d79=\"],\n  [\"656\u00d72=\", \"577\u00d76=\"],\n  [\"554\u00d73=\", \"944\u00d72=\"],\n  [\"120\u00d74=\", \"609\u00d75=\"],\n  [\"589\u00d77=\", \"144\u00d73=\"],\n  [\"930\u00d74=\", \"199\u00d73=\"],\n  [\"856\u00d72=\", \"612\u00d75=\"],\n  [\"457\u00d74=\", \"138\u00d76=\"],\n  [\"514\u00d78=\", \"891\u00d77=\"],\n  [\"613\u00d76=\", \"759\u00d79=\"],\n  [\"187\u00d74=\", \"239\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const item of found.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression in the\n# document with its new value, per the commit's regenerated numbers.\n# Every \"old\" string is unique in the document, so a MatchCase\n# Find/Replace-All for each pair unambiguously targets the single cell\n# that holds it.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"195\u00d77=\", \"761\u00d75=\"),\n  @(\"552\u00d78=\", \"230\u00d73=\"),\n  @(\"709\u00d79=\", \"367\u00d75=\"),\n  @(\"248\u00d76=\", \"844\u00d79=\"),\n  @(\"234\u00d77=\", \"418\u00d72=\"),\n  @(\"305\u00d79=\", \"181\u00d75=\"),\n  @(\"260\u00d74=\", \"261\u00d73=\"),\n  @(\"409\u00d73=\", \"423\u00d76=\"),\n  @(\"314\u00d77=\", \"721\u00d73=\"),\n  @(\"149\u00d74=\", \"935\u00d72=\"),\n  @(\"243\u00d75=\", \"698\u00d72=\"),\n  @(\"386\u00d72=\", \"484\u00d78=\"),\n  @(\"182\u00d77=\", \"513\u00d72=\"),\n  @(\"976\u00d75=\", \"814\u00d74=\"),\n  @(\"778\u00d75=\", \"620\u00d79=\"),\n  @(\"656\u00d72=\", \"577\u00d76=\"),\n  @(\"554\u00d73=\", \"944\u00d72=\"),\n  @(\"120\u00d74=\", \"609\u00d75=\"),\n  @(\"589\u00d77=\", \"144\u00d73=\"),\n  @(\"930\u00d74=\", \"199\u00d73=\"),\n  @(\"856\u00d72=\", \"612\u00d75=\"),\n  @(\"457\u00d74=\", \"138\u00d76=\"),\n  @(\"514\u00d78=\", \"891\u00d77=\"),\n  @(\"613\u00d76=\", \"759\u00d79=\"),\n  @(\"187\u00d74=\", \"239\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, \"wdFindContinue\", $false, $find.Replacement.Text, \"wdReplaceAll\")\n}\n\nWrite-Output \"replacements complete\"\n"}
